# Rename data sheets from their hash-based names to simple sequential
# index numbers ("2".."14"), and switch the "Index" sheet's first column
# from a text reference of that hash to the numeric sheet index.

$wb = $excel.ActiveWorkbook

# Work from the 2nd worksheet (index 2) through the 14th (index 14) -
# rename each to its 1-based worksheet position (as a string).
for ($i = 2; $i -le 14; $i++) {
    $sheet = $wb.Worksheets.Item($i)
    $sheet.Name = "$i"
}

# Update the "Index" sheet: column A used to hold the sheet-name hash
# (a shared string); it now holds the plain numeric sheet index.
$idx = $wb.Worksheets.Item("Index")
for ($r = 2; $r -le 14; $r++) {
    $idx.Cells.Item($r, 1).Value = $r
}
